$wb = $excel.ActiveWorkbook

# Rename the existing sheet to "Invoice"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Invoice"

# Add a new sheet "Customer" right after "Invoice"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Customer"

# Populate the Customer sheet
$ws2.Range("A1").Value = "customer_id"
$ws2.Range("B1").Value = "customer_name"
$ws2.Range("A2").Value = 1439
$ws2.Range("B2").Value = "Lester Chalmers"
